# Fix layout issue with group ID on contestant cards
# - Remove the "Standbys" sheet: the one standby contestant (Kathleen
#   Reynolds) is now promoted straight to a confirmed Seat Assignment.
# - Insert the promoted contestant as a new row in "Seat Assignments",
#   pushing the existing row down.
# - Add an explicit (empty) H4 cell on "Contestants" so the row is fully
#   populated in column order.

$wb = $excel.ActiveWorkbook

# --- Contestants: add an explicit empty cell at H4 -------------------------
# (Rating is blank for this contestant; touch the cell's format so a real
#  -- if empty -- cell entry is materialized at H4 instead of leaving a gap
#  in the row, matching the other blank-but-present cells in this sheet.)
$wsContestants = $wb.Worksheets.Item("Contestants")
$wsContestants.Range("H4").Font.Size = 12

# --- Seat Assignments: insert the (former) standby as a new row 3 ----------
$wsSeats = $wb.Worksheets.Item("Seat Assignments")
$wsSeats.Rows.Item(3).Insert()
$wsSeats.Range("A3").Value = "479fcb75-5c6b-42a9-a757-6e4ae8cef47e"
$wsSeats.Range("B3").Value = "e432f0fe-1383-44a2-990c-5f787da5008a"
$wsSeats.Range("C3").Value = "d698b1de-6641-45c6-aa63-f577d2b634bb"
$wsSeats.Range("D3").Value = 1
$wsSeats.Range("E3").Value = "A4"

# --- Standbys: no longer needed, remove the whole sheet --------------------
$wsStandbys = $wb.Worksheets.Item("Standbys")
$wsStandbys.Delete()
